# Add data for 2021-11-20
# (source data actually reflects "through 11-12", per sheet name / header text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet/tab to reflect the new "through" date.
$ws.Name = "Through 2021-11-12"

# --- Row 12 (October) ---
# Only the 2021 arrest_made/no_arrest_made/arrest_rate triple (T:V) changes.
$ws.Range("U12").Value = 190
$ws.Range("V12").Value = 0.0206

# --- Row 13 (November, through-date label + counts) ---
$ws.Range("A13").Value = "November (through 11-12)"

$ws.Range("C13").Value = 14

$ws.Range("F13").Value = 27
$ws.Range("G13").Value = 0.0357

$ws.Range("I13").Value = 48
$ws.Range("J13").Value = 0.0204

$ws.Range("L13").Value = 21
$ws.Range("M13").Value = 0.1923

$ws.Range("N13").Value = 4
$ws.Range("O13").Value = 17
$ws.Range("P13").Value = 0.1905

$ws.Range("R13").Value = 75
$ws.Range("S13").Value = 0.026

$ws.Range("U13").Value = 82
$ws.Range("V13").Value = 0.012

# --- Row 14 (Total) ---
$ws.Range("C14").Value = 240
$ws.Range("D14").Value = 0.1176

$ws.Range("F14").Value = 461
$ws.Range("G14").Value = 0.1031

$ws.Range("I14").Value = 697
$ws.Range("J14").Value = 0.0817

$ws.Range("L14").Value = 570
$ws.Range("M14").Value = 0.1108

$ws.Range("N14").Value = 52
$ws.Range("O14").Value = 451
$ws.Range("P14").Value = 0.1034

$ws.Range("R14").Value = 1078
$ws.Range("S14").Value = 0.0494

$ws.Range("U14").Value = 1436
$ws.Range("V14").Value = 0.0584
